$wb = $excel.ActiveWorkbook

# ---- Sheet: ALC ----
$ws = $wb.Worksheets.Item("ALC")
# Row 40
$ws.Range("H40").Value = 2799
$ws.Range("J40").Value = 2799
$ws.Range("L40").Value = 2799
$ws.Range("N40").Value = -3149
# Row 62
$ws.Range("H62").Value = 4508.4546
$ws.Range("I62").Value = 4470.4
$ws.Range("J62").Value = 4889
$ws.Range("K62").Value = 4470.4
$ws.Range("L62").Value = 4889
$ws.Range("M62").Value = -3846.4
$ws.Range("N62").Value = -6137
# Row 65
$ws.Range("H65").Value = 4508.4546
$ws.Range("I65").Value = 4470.4
$ws.Range("J65").Value = 4889
$ws.Range("K65").Value = 22352
$ws.Range("L65").Value = 24445
$ws.Range("M65").Value = -19232
$ws.Range("N65").Value = -30685
# Row 70
$ws.Range("H70").Value = 2147.5
$ws.Range("I70").Value = 1845
$ws.Range("K70").Value = 5535
$ws.Range("M70").Value = -5265
# Row 73
$ws.Range("H73").Value = 2147.5
$ws.Range("I73").Value = 1845
$ws.Range("K73").Value = 5535
$ws.Range("M73").Value = -4599
# Row 92
$ws.Range("H92").Value = 1532
$ws.Range("I92").Value = 1415
$ws.Range("K92").Value = 1415
$ws.Range("M92").Value = -167
# Row 132
$ws.Range("H132").Value = 28573528
$ws.Range("I132").Value = 32260196
$ws.Range("K132").Value = 96780588
$ws.Range("M132").Value = -96778058
# Row 135
$ws.Range("H135").Value = 5283.48
$ws.Range("I135").Value = 4512.5557
$ws.Range("J135").Value = 7265.857
$ws.Range("K135").Value = 40613.0013
$ws.Range("L135").Value = 65392.713
$ws.Range("M135").Value = -38078.0013
$ws.Range("N135").Value = -70462.713
# Row 137
$ws.Range("H137").Value = 1735.0667
$ws.Range("I137").Value = 1716.1428
$ws.Range("K137").Value = 5148.428400000001
$ws.Range("M137").Value = -2598.428400000001
# Row 138
$ws.Range("H138").Value = 3583.1086
$ws.Range("J138").Value = 6693.45
$ws.Range("L138").Value = 20080.35
$ws.Range("N138").Value = -30360.35

# ---- Sheet: ARM ----
$ws = $wb.Worksheets.Item("ARM")
# Row 32
$ws.Range("H32").Value = 7220.758
$ws.Range("I32").Value = 5408.478
$ws.Range("J32").Value = 11389
$ws.Range("K32").Value = 5408.478
$ws.Range("L32").Value = 11389
$ws.Range("M32").Value = -5121.478
$ws.Range("N32").Value = -11963
# Row 61
$ws.Range("H61").Value = 5639.9443
$ws.Range("I61").Value = 2720.3333
$ws.Range("K61").Value = 2720.3333
$ws.Range("M61").Value = -2508.3333
# Row 74
$ws.Range("H74").Value = 2396.5588
$ws.Range("I74").Value = 2149.8
$ws.Range("K74").Value = 2149.8
$ws.Range("M74").Value = -1275.8
# Row 77
$ws.Range("H77").Value = 2396.5588
$ws.Range("I77").Value = 2149.8
$ws.Range("K77").Value = 10749
$ws.Range("M77").Value = -6381
# Row 136
$ws.Range("H136").Value = 5639.9443
$ws.Range("I136").Value = 2720.3333
$ws.Range("K136").Value = 8160.999899999999
$ws.Range("M136").Value = -5610.999899999999

# ---- Sheet: BSM ----
$ws = $wb.Worksheets.Item("BSM")
# Row 20
$ws.Range("H20").Value = 43480892
$ws.Range("I20").Value = 62502264
$ws.Range("J20").Value = 3468.4285
$ws.Range("K20").Value = 62502264
$ws.Range("L20").Value = 3468.4285
$ws.Range("M20").Value = -62502017
$ws.Range("N20").Value = -3962.4285
# Row 134
$ws.Range("H134").Value = 16130883
$ws.Range("I134").Value = 17243288
$ws.Range("K134").Value = 51729864
$ws.Range("M134").Value = -51727329
# Row 140
$ws.Range("H140").Value = 89909
$ws.Range("J140").Value = 89909
$ws.Range("L140").Value = 89909
$ws.Range("N140").Value = -100269

# ---- Sheet: CRP ----
$ws = $wb.Worksheets.Item("CRP")
# Row 31
$ws.Range("H31").Value = 2986
$ws.Range("I31").Value = 3259.75
$ws.Range("K31").Value = 3259.75
$ws.Range("M31").Value = -2964.75
# Row 34
$ws.Range("H34").Value = 2986
$ws.Range("I34").Value = 3259.75
$ws.Range("K34").Value = 3259.75
$ws.Range("M34").Value = -3057.75
# Row 58
$ws.Range("H58").Value = 4686.5
$ws.Range("I58").Value = 4748.9473
$ws.Range("K58").Value = 4748.9473
$ws.Range("M58").Value = -4545.9473
# Row 74
$ws.Range("H74").Value = 72899.8
$ws.Range("J74").Value = 72899.8
$ws.Range("L74").Value = 72899.8
$ws.Range("N74").Value = -74647.8
# Row 77
$ws.Range("H77").Value = 72899.8
$ws.Range("J77").Value = 72899.8
$ws.Range("L77").Value = 218699.4
$ws.Range("N77").Value = -227435.4
# Row 86
$ws.Range("H86").Value = 76926740
$ws.Range("I86").Value = 90912360
$ws.Range("K86").Value = 90912360
$ws.Range("M86").Value = -90911237
# Row 89
$ws.Range("H89").Value = 76926740
$ws.Range("I89").Value = 90912360
$ws.Range("K89").Value = 454561800
$ws.Range("M89").Value = -454556184
# Row 94
$ws.Range("H94").Value = 1603.6111
$ws.Range("I94").Value = 1577.9
$ws.Range("J94").Value = 1635.75
$ws.Range("K94").Value = 1577.9
$ws.Range("L94").Value = 1635.75
$ws.Range("M94").Value = -1126.9
$ws.Range("N94").Value = -2537.75
# Row 105
$ws.Range("H105").Value = 12620.625
$ws.Range("I105").Value = 13328.667
$ws.Range("K105").Value = 13328.667
$ws.Range("M105").Value = -11581.667
# Row 134
$ws.Range("H134").Value = 1512.5333
$ws.Range("I134").Value = 1349.0834
$ws.Range("K134").Value = 4047.2502
$ws.Range("M134").Value = -1512.2502
# Row 136
$ws.Range("H136").Value = 4686.5
$ws.Range("I136").Value = 4748.9473
$ws.Range("K136").Value = 14246.8419
$ws.Range("M136").Value = -11696.8419

# ---- Sheet: CUL ----
$ws = $wb.Worksheets.Item("CUL")
# Row 64
$ws.Range("H64").Value = 3971.5454
$ws.Range("I64").Value = 1711.125
$ws.Range("J64").Value = 9999.333000000001
$ws.Range("K64").Value = 5133.375
$ws.Range("L64").Value = 29997.999
$ws.Range("M64").Value = -4863.375
$ws.Range("N64").Value = -30537.999
# Row 67
$ws.Range("H67").Value = 3971.5454
$ws.Range("I67").Value = 1711.125
$ws.Range("J67").Value = 9999.333000000001
$ws.Range("K67").Value = 5133.375
$ws.Range("L67").Value = 29997.999
$ws.Range("M67").Value = -4197.375
$ws.Range("N67").Value = -31869.999
# Row 121
$ws.Range("H121").Value = 12338554
$ws.Range("I121").Value = 501500
$ws.Range("J121").Value = 15720569
$ws.Range("K121").Value = 1504500
$ws.Range("L121").Value = 47161707
$ws.Range("M121").Value = -1503190
$ws.Range("N121").Value = -47164327
# Row 123
$ws.Range("H123").Value = 4444
$ws.Range("I123").Value = 4444
$ws.Range("K123").Value = 13332
$ws.Range("M123").Value = -10882
# Row 132
$ws.Range("H132").Value = 2110.9092
$ws.Range("J132").Value = 2625
$ws.Range("L132").Value = 23625
$ws.Range("N132").Value = -28685

# ---- Sheet: GSM ----
$ws = $wb.Worksheets.Item("GSM")
# Row 33
$ws.Range("H33").Value = 34832.918
$ws.Range("J33").Value = 34832.918
$ws.Range("L33").Value = 34832.918
$ws.Range("N33").Value = -35336.918
# Row 70
$ws.Range("H70").Value = 28949.852
$ws.Range("I70").Value = 49984.75
$ws.Range("J70").Value = 12121.934
$ws.Range("K70").Value = 49984.75
$ws.Range("L70").Value = 12121.934
$ws.Range("M70").Value = -49714.75
$ws.Range("N70").Value = -12661.934
# Row 73
$ws.Range("H73").Value = 28949.852
$ws.Range("I73").Value = 49984.75
$ws.Range("J73").Value = 12121.934
$ws.Range("K73").Value = 49984.75
$ws.Range("L73").Value = 12121.934
$ws.Range("M73").Value = -49048.75
$ws.Range("N73").Value = -13993.934
# Row 88
$ws.Range("H88").Value = 0
$ws.Range("J88").Value = 0
$ws.Range("L88").Value = 0
$ws.Range("N88").ClearContents()
# Row 91
$ws.Range("H91").Value = 0
$ws.Range("J91").Value = 0
$ws.Range("L91").Value = 0
$ws.Range("N91").ClearContents()
# Row 102
$ws.Range("H102").Value = 1353.8823
$ws.Range("I102").Value = 1353.8823
$ws.Range("K102").Value = 1353.8823
$ws.Range("M102").Value = 268.1177
# Row 126
$ws.Range("H126").Value = 6679.2856
$ws.Range("I126").Value = 7059.3335
$ws.Range("J126").Value = 4399
$ws.Range("K126").Value = 21178.0005
$ws.Range("L126").Value = 13197
$ws.Range("M126").Value = -18708.0005
$ws.Range("N126").Value = -18137

# ---- Sheet: LTW ----
$ws = $wb.Worksheets.Item("LTW")
# Row 38
$ws.Range("H38").Value = 37816.332
$ws.Range("I38").Value = 37799
$ws.Range("J38").Value = 37825
$ws.Range("K38").Value = 37799
$ws.Range("L38").Value = 37825
$ws.Range("M38").Value = -37389
$ws.Range("N38").Value = -38645
# Row 61
$ws.Range("H61").Value = 1022.7059
$ws.Range("I61").Value = 886.875
$ws.Range("K61").Value = 886.875
$ws.Range("M61").Value = -684.875
# Row 113
$ws.Range("H113").Value = 1022.7059
$ws.Range("I113").Value = 886.875
$ws.Range("K113").Value = 886.875
$ws.Range("M113").Value = 1283.125
# Row 122
$ws.Range("H122").Value = 8231.473
$ws.Range("I122").Value = 8165.7
$ws.Range("K122").Value = 24497.1
$ws.Range("M122").Value = -22047.1
# Row 132
$ws.Range("H132").Value = 5060.2144
$ws.Range("I132").Value = 5077.788
$ws.Range("K132").Value = 15233.364
$ws.Range("M132").Value = -12703.364
# Row 136
$ws.Range("H136").Value = 6730.5654
$ws.Range("I136").Value = 5407.1113
$ws.Range("K136").Value = 16221.3339
$ws.Range("M136").Value = -13671.3339

# ---- Sheet: WVR ----
$ws = $wb.Worksheets.Item("WVR")
# Row 62
$ws.Range("H62").Value = 35728484
$ws.Range("J62").Value = 7255
$ws.Range("L62").Value = 7255
$ws.Range("N62").Value = -8503
# Row 65
$ws.Range("H65").Value = 35728484
$ws.Range("J65").Value = 7255
$ws.Range("L65").Value = 36275
$ws.Range("N65").Value = -42515
# Row 136
$ws.Range("H136").Value = 4725.6895
$ws.Range("I136").Value = 4088.652
$ws.Range("K136").Value = 12265.956
$ws.Range("M136").Value = -9715.956
# Row 138
$ws.Range("H138").Value = 60000
$ws.Range("J138").Value = 60000
$ws.Range("L138").Value = 60000
$ws.Range("N138").Value = -70280
